$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (post-edit) order of player rows 2-19, as (Name, Position, Team)
$data = @(
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("Jabari Smith Jr.", "PF,C", "Houston Rockets"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("P.J. Washington", "PF", "Dallas Mavericks"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Dereck Lively II", "C", "Dallas Mavericks"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Norman Powell", "SG,SF", "LA Clippers")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
